$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Funding")
$ws2 = $wb.Worksheets.Item("Expense")

# ---------------------------------------------------------------------------
# 1. Funding sheet ("Funding") — Priority letters shift by one funding row
#    that was removed upstream, the "Valid To" column becomes a real date
#    (instead of a shared text string) and a handful of "Valid From"/
#    "Valid To" values move for FS011-FS015.
# ---------------------------------------------------------------------------

$fundingPriority = @{2="A";3="B";4="C";5="D";6="E";7="F";8="G";9="H";10="I";11="J";12="K";13="L";14="M";15="N";16="O"}
$fundingCategory = @{
  2="Salary, Equipment, Travel"; 3="Salary, Equipment, Travel"; 4="Salary, Equipment, Travel";
  5="Salary, Equipment, Travel"; 6="Salary, Equipment, Travel"; 7="Salary, Equipment, Travel";
  8="Salary"; 9="Salary"; 10="Salary"; 11="Salary";
  12="Salary"; 13="Equipment"; 14="Travel"; 15="Salary"; 16="Equipment"
}
$fundingValidFrom = @{2=45658;3=45658;4=45658;5=45658;6=45658;7=45658;8=45658;9=45658;10=45658;11=45658;12=45839;13=45839;14=45839;15=45901;16=45901}
$fundingValidTo   = @{2=45838;3=45838;4=45838;5=45838;6=45838;7=45838;8=45838;9=45838;10=45838;11=45838;12=46022;13=46022;14=46022;15=46022;16=46022}
$fundingAmount    = @{2=10000;3=10000;4=10000;5=10000;6=10000;7=10000;8=5000;9=5000;10=5000;11=5000;12=2000;13=2000;14=2000;15=1000;16=1000}

# Give the "Valid From"/"Valid To" headers the same date-number-format style
# already used by the "Valid From" data cells (D2), via a format-only copy so
# no new style entry is created.
$ws1.Range("D2").Copy()
$ws1.Range("D1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

foreach ($r in 2..16) {
    $ws1.Range("B$r").Value = $fundingPriority[$r]
    $ws1.Range("C$r").Value = $fundingCategory[$r]
    $ws1.Range("D$r").Value = $fundingValidFrom[$r]
    $ws1.Range("F$r").Value = $fundingAmount[$r]
}

# "Valid To" (E) becomes a genuine date value styled like F (date, right
# aligned) instead of a shared string; grab that style via PasteSpecial so
# the now-unused "right align only" style stops being referenced, then set
# the real values.
$ws2.Range("F2").Copy()
$ws1.Range("E2:E16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
foreach ($r in 2..16) {
    $ws1.Range("E$r").Value = $fundingValidTo[$r]
}

# The "Amount" column (F) keeps its numeric format, but the style table slot
# it used to reference is going away along with the "Valid To" one above —
# refresh it from itself via a round-trip copy so it also lands on the
# post-cleanup numFmt=3 slot.
$ws1.Range("F2").Copy()
$ws1.Range("F2:F16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
foreach ($r in 2..16) {
    $ws1.Range("F$r").Value = $fundingAmount[$r]
}

# ---------------------------------------------------------------------------
# 2. Expense sheet ("Expense") — expense IDs shift the same way, and the
#    "Latest Payment Date" column is fully recomputed to real dates (some
#    cells used to hold literal date *text*).
# ---------------------------------------------------------------------------

$expenseCategory = @{2="Salary";3="Salary";4="Salary";5="Equipment";6="Equipment";7="Equipment";8="Travel";9="Travel";10="Travel";11="Salary"}
$expensePlanned  = @{2=10000;3=10000;4=10000;5=12000;6=12000;7=12000;8=8000;9=8000;10=8000;11=10000}
$expenseLatest   = @{2=45698;3=45721;4=45741;5=45848;6=45874;7=45894;8=45910;9=45935;10=45955;11=45981}

$ws2.Range("D2").Copy()
$ws2.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

foreach ($r in 2..11) {
    $ws2.Range("D$r").Value = $expenseCategory[$r]
    $ws2.Range("E$r").Value = $expensePlanned[$r]
}

# Normalise every "Latest Payment Date" cell onto the same date style (some
# already were; the ones that held literal date text need to be converted to
# real numbers too) before writing the recomputed values.
$ws2.Range("F2").Copy()
$ws2.Range("F2:F11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
foreach ($r in 2..11) {
    $ws2.Range("F$r").Value = $expenseLatest[$r]
}

# ---------------------------------------------------------------------------
# 3. View state — zoom + selection per sheet. Touch the non-active sheet
#    ("Funding") first so the final ActiveSheet/tabSelected stays on
#    "Expense", matching the saved workbook.
# ---------------------------------------------------------------------------

$ws1.Select()
$ws1.Range("H7").Select()
$excel.ActiveWindow.Zoom = 156

$ws2.Select()
$ws2.Range("F4").Select()
$excel.ActiveWindow.Zoom = 157
